# fix: Reconocer si un query es de siniestros, primas, o expuestos
#
# The sheet "add_pe_Amparos" was generalized/renamed to "add_e_Amparos"
# (no longer tied to "primas y expuestos" only). Renaming via the
# Worksheets collection also keeps every reference (defined names, the
# sheet's own AutoFilter defined name, etc.) in sync automatically, the
# same way Excel's UI rename does.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("add_pe_Amparos")
$ws.Name = "add_e_Amparos"

# Reflect that this sheet was the one being worked on when the file was
# saved: make it the active tab with cell J21 selected.
$ws.Activate() | Out-Null
$ws.Range("J21").Select() | Out-Null
